$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: was "収支入力" admin function row -> becomes the admin "edit user" row
$ws.Range("B8").Value = "ユーザーを編集する"
$ws.Range("C8").Value = "管理者がユーザー情報を編集する"
$ws.Range("D8").Value = "ユーザーID・名前・種別"
$ws.Range("E8").Value = "管理者"

# Row 9: becomes the admin "delete user" row
$ws.Range("B9").Value = "ユーザーを削除する"
$ws.Range("C9").Value = "管理者がユーザー情報を削除する"
$ws.Range("D9").Value = "ユーザーID"
$ws.Range("E9").Value = "管理者"

# Row 10: shifts to "収支入力" (income/expense entry)
$ws.Range("B10").Value = "収支入力"
$ws.Range("C10").Value = "収入・支出データを入力"
$ws.Range("E10").Value = "管理者・利用者"

# Row 11: shifts to "収支編集" (income/expense edit)
$ws.Range("B11").Value = "収支編集"
$ws.Range("C11").Value = "登録済データを修正"
$ws.Range("D11").Value = "日付・金額・カテゴリ・メモ"
$ws.Range("E11").Value = "管理者・利用者"

# Row 12: shifts to "収支削除" (income/expense delete); F12 no longer used
$ws.Range("B12").Value = "収支削除"
$ws.Range("C12").Value = "収支データを削除"
$ws.Range("D12").Value = "日付・金額・カテゴリ・メモ"
$ws.Range("E12").Value = "管理者・利用者"
$ws.Range("F12").ClearContents()

# Row 13: shifts to "収支一覧表示" (income/expense list view)
$ws.Range("B13").Value = "収支一覧表示"
$ws.Range("C13").Value = "登録済データを一覧表示"
$ws.Range("D13").Value = "日付・金額・カテゴリ"

# Row 14: shifts to "グラフ表示" (graph display); gains F14
$ws.Range("B14").Value = "グラフ表示"
$ws.Range("C14").Value = "月別・カテゴリ別で可視化"
$ws.Range("D14").Value = "収支情報"
$ws.Range("E14").Value = "管理者・利用者"
$ws.Range("F14").Value = "チャレンジ要素"

# Row 15: shifts to "目標金額設定" (goal amount setting)
$ws.Range("B15").Value = "目標金額設定"
$ws.Range("C15").Value = "月ごとの支出目標を設定"
$ws.Range("D15").Value = "月・金額"

# Row 16 (new): "達成状況表示" (achievement status display)
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "達成状況表示"
$ws.Range("C16").Value = "目標に対する進捗を表示"
$ws.Range("D16").Value = "目標金額・実績"
$ws.Range("E16").Value = "管理者・利用者"

# Row 17 (new): "進捗判定" (progress determination)
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "進捗判定"
$ws.Range("C17").Value = "目標達成かどうか判定"
$ws.Range("D17").Value = "目標金額・実績"
$ws.Range("E17").Value = "管理者・利用者"

# Update selection to match the final saved view state
$ws.Range("D9").Select()
